# Day 1 Array Refresher PPT
# Insert two new paragraphs at the very beginning of the document body:
#   1) A bold "SOLVE 2, 4 and 5" paragraph with a hanging indent.
#   2) An empty ListParagraph-styled paragraph.

$d = $word.ActiveDocument

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:ind w:left="720" w:hanging="360"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>SOLVE 2, 4 and 5</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
            </w:pPr>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

# Collapsed range at the very start of the document; InsertXML there drops
# the two new paragraphs in before all existing content without disturbing
# the formatting/numbering of the paragraph that used to be first.
$target = $d.Range(0, 0)
$target.InsertXML($xml)
